$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.529888943768071
$ws.Range("D2").Value = 0.001752524866425631
$ws.Range("E2").Value = 0.7367889436378476
$ws.Range("F2").Value = 0.7518689970549275
$ws.Range("G2").Value = 0.6439444612340708
$ws.Range("H2").Value = 0.6092161822974447
$ws.Range("I2").Value = 0.7325154218006062
$ws.Range("L2").Value = 0.6142176853802255
$ws.Range("B3").Value = 1.391734653961521
$ws.Range("D3").Value = 0.001956490648494569
$ws.Range("E3").Value = 0.6750170817489192
$ws.Range("F3").Value = 0.7132351469910247
$ws.Range("G3").Value = 0.6000875643137249
$ws.Range("H3").Value = 0.5962483534734986
$ws.Range("I3").Value = 0.7397454770155463
$ws.Range("L3").Value = 0.5452547963602399
$ws.Range("B4").Value = 1.306782738698587
$ws.Range("D4").Value = 0.00209386072557205
$ws.Range("E4").Value = 0.6369834076335081
$ws.Range("F4").Value = 0.6903673921852942
$ws.Range("G4").Value = 0.5739508704097034
$ws.Range("H4").Value = 0.5889385751365808
$ws.Range("I4").Value = 0.7450211596285357
$ws.Range("L4").Value = 0.5028710331772572
$ws.Range("B5").Value = 1.272134359397967
$ws.Range("D5").Value = 0.00215283144169387
$ws.Range("E5").Value = 0.6214591680858064
$ws.Range("F5").Value = 0.6812602609967229
$ws.Range("G5").Value = 0.5634954389955453
$ws.Range("H5").Value = 0.5861223338666264
$ws.Range("I5").Value = 0.7473810570008723
$ws.Range("L5").Value = 0.4855895503614533
$ws.Range("B6").Value = 1.266379270367395
$ws.Range("D6").Value = 0.002162802405517361
$ws.Range("E6").Value = 0.6188798966917233
$ws.Range("F6").Value = 0.6797607260581913
$ws.Range("G6").Value = 0.5617710250336927
$ws.Range("H6").Value = 0.5856644723756403
$ws.Range("I6").Value = 0.7477855919617369
$ws.Range("L6").Value = 0.4827193955854057
$ws.Range("B7").Value = 1.306315577044359
$ws.Range("D7").Value = 0.002094643992916367
$ws.Range("E7").Value = 0.6367741425993643
$ws.Range("F7").Value = 0.6902437171416409
$ws.Range("G7").Value = 0.5738090780173479
$ws.Range("H7").Value = 0.5888999381834594
$ws.Range("I7").Value = 0.7450521360899387
$ws.Range("L7").Value = 0.5026380080422257
$ws.Range("B8").Value = 1.482280083169712
$ws.Range("D8").Value = 0.001820297619985478
$ws.Range("E8").Value = 0.7155126688945757
$ws.Range("F8").Value = 0.738369075601824
$ws.Range("G8").Value = 0.6286561333906491
$ws.Range("H8").Value = 0.6046085159035215
$ws.Range("I8").Value = 0.7348345577386581
$ws.Range("L8").Value = 0.5904478336959755
$ws.Range("B9").Value = 1.826308372257017
$ws.Range("D9").Value = 0.001381282631902625
$ws.Range("E9").Value = 0.8690318839206554
$ws.Range("F9").Value = 0.839655556556437
$ws.Range("G9").Value = 0.7426593001476363
$ws.Range("H9").Value = 0.6406620857524103
$ws.Range("I9").Value = 0.7214502745963074
$ws.Range("L9").Value = 0.7623165481452077
$ws.Range("B10").Value = 2.078396669518895
$ws.Range("D10").Value = 0.001122802391431588
$ws.Range("E10").Value = 0.9812273147321378
$ws.Range("F10").Value = 0.9184870555936158
$ws.Range("G10").Value = 0.8305883093823638
$ws.Range("H10").Value = 0.6704536846989413
$ws.Range("I10").Value = 0.7156973627234393
$ws.Range("L10").Value = 0.8883953150881041
$ws.Range("B11").Value = 2.192927737255445
$ws.Range("D11").Value = 0.001019936876050842
$ws.Range("E11").Value = 1.032128261231719
$ws.Range("F11").Value = 0.955352868216778
$ws.Range("G11").Value = 0.8715476809360041
$ws.Range("H11").Value = 0.6847465092708376
$ws.Range("I11").Value = 0.7139722454856212
$ws.Range("L11").Value = 0.9457123783432735
$ws.Range("B12").Value = 2.236275874300304
$ws.Range("D12").Value = 0.0009831691215866911
$ws.Range("E12").Value = 1.05138229694029
$ws.Range("F12").Value = 0.9694609795060813
$ws.Range("G12").Value = 0.8872001123072266
$ws.Range("H12").Value = 0.6902670984166548
$ws.Range("I12").Value = 0.7134477449650873
$ws.Range("L12").Value = 0.9674115178373768
$ws.Range("B13").Value = 2.226941090824255
$ws.Range("D13").Value = 0.0009909894258952612
$ws.Range("E13").Value = 1.047236552426085
$ws.Range("F13").Value = 0.9664159175005551
$ws.Range("G13").Value = 0.883822702027544
$ws.Range("H13").Value = 0.6890733005585901
$ws.Range("I13").Value = 0.7135549703799953
$ws.Range("L13").Value = 0.962738472884098
$ws.Range("B14").Value = 2.196494469275365
$ws.Range("D14").Value = 0.001016867780838648
$ws.Range("E14").Value = 1.033712731535388
$ws.Range("F14").Value = 0.9565105692293514
$ws.Range("G14").Value = 0.8728325470850677
$ws.Range("H14").Value = 0.6851985128863021
$ws.Range("I14").Value = 0.7139265108964636
$ws.Range("L14").Value = 0.9474976932123411
$ws.Range("B15").Value = 2.177842090181684
$ws.Range("D15").Value = 0.00103300564893205
$ws.Range("E15").Value = 1.025426221478966
$ws.Range("F15").Value = 0.9504626079000218
$ws.Range("G15").Value = 0.8661193712553654
$ws.Range("H15").Value = 0.6828392368067
$ws.Range("I15").Value = 0.7141708745414448
$ws.Range("L15").Value = 0.9381615394746632
$ws.Range("B16").Value = 2.070908751260902
$ws.Range("D16").Value = 0.001129826919134214
$ws.Range("E16").Value = 0.9778979431024766
$ws.Range("F16").Value = 0.9160982932412765
$ws.Range("G16").Value = 0.8279311736621651
$ws.Range("H16").Value = 0.6695346625245122
$ws.Range("I16").Value = 0.7158280936275929
$ws.Range("L16").Value = 0.8846487437281496
$ws.Range("B17").Value = 2.00527040399453
$ws.Range("D17").Value = 0.001193041973963638
$ws.Range("E17").Value = 0.9487047788308445
$ws.Range("F17").Value = 0.8952766447480798
$ws.Range("G17").Value = 0.8047526406831764
$ws.Range("H17").Value = 0.6615636144032919
$ws.Range("I17").Value = 0.717073537794306
$ws.Range("L17").Value = 0.8518107318939485
$ws.Range("B18").Value = 1.967503380357073
$ws.Range("D18").Value = 0.001230782809858777
$ws.Range("E18").Value = 0.9319007940478059
$ws.Range("F18").Value = 0.8833949756088515
$ws.Range("G18").Value = 0.7915110850485121
$ws.Range("H18").Value = 0.6570484712625557
$ws.Range("I18").Value = 0.7178738033937151
$ws.Range("L18").Value = 0.8329197282417624
$ws.Range("B19").Value = 1.95471382647122
$ws.Range("D19").Value = 0.001243796353956661
$ws.Range("E19").Value = 0.9262090868386395
$ws.Range("F19").Value = 0.8793881580585321
$ws.Range("G19").Value = 0.7870430806382558
$ws.Range("H19").Value = 0.6555316216726226
$ws.Range("I19").Value = 0.7181591576455375
$ws.Range("L19").Value = 0.8265229763188131
$ws.Range("B20").Value = 2.012259134294766
$ws.Range("D20").Value = 0.001186169178449159
$ws.Range("E20").Value = 0.9518137803201796
$ws.Range("F20").Value = 0.8974833516116263
$ws.Range("G20").Value = 0.8072106791025817
$ws.Range("H20").Value = 0.6624049310283908
$ws.Range("I20").Value = 0.7169322692296234
$ws.Range("L20").Value = 0.8553067533769934
$ws.Range("B21").Value = 2.205437996092257
$ws.Range("D21").Value = 0.001009206823348041
$ws.Range("E21").Value = 1.037685587038482
$ws.Range("F21").Value = 0.9594159715402668
$ws.Range("G21").Value = 0.8760567362893994
$ws.Range("H21").Value = 0.686333681403255
$ws.Range("I21").Value = 0.7138138816564137
$ws.Range("L21").Value = 0.9519744335428584
$ws.Range("B22").Value = 2.331560993789196
$ws.Range("D22").Value = 0.0009063211254769499
$ws.Range("E22").Value = 1.09368450701399
$ws.Range("F22").Value = 1.000755735055037
$ws.Range("G22").Value = 0.9218811964900908
$ws.Range("H22").Value = 0.7026040282694339
$ws.Range("I22").Value = 0.7125266224499001
$ws.Range("L22").Value = 1.015119658535127
$ws.Range("B23").Value = 2.264259175755228
$ws.Range("D23").Value = 0.0009600418270925282
$ws.Range("E23").Value = 1.063808532157026
$ws.Range("F23").Value = 0.9786118375208304
$ws.Range("G23").Value = 0.8973465967743834
$ws.Range("H23").Value = 0.6938618781725268
$ws.Range("I23").Value = 0.7131447857711493
$ws.Range("L23").Value = 0.9814209459457288
$ws.Range("B24").Value = 2.009099624726502
$ws.Range("D24").Value = 0.001189272020745591
$ws.Range("E24").Value = 0.9504082643316139
$ws.Range("F24").Value = 0.8964854225549601
$ws.Range("G24").Value = 0.8060991383050862
$ws.Range("H24").Value = 0.6620243616070525
$ws.Range("I24").Value = 0.7169958743872158
$ws.Range("L24").Value = 0.853726239397588
$ws.Range("B25").Value = 1.73335470640734
$ws.Range("D25").Value = 0.001489069567629731
$ws.Range("E25").Value = 0.8276015395192218
$ws.Range("F25").Value = 0.8114940263579626
$ws.Range("G25").Value = 0.7111034090098656
$ws.Range("H25").Value = 0.6303360915817109
$ws.Range("I25").Value = 0.7243563959884867
$ws.Range("L25").Value = 0.7158559362988797
